# Weekly update: a new price observation row is inserted as the new first
# data row (row 31) for "Comercializadora del Agro de Limarí - Poroto
# granado", pushing all subsequent rows down by one. The new row carries
# the same categorical/fixed fields (mercado, región, codreg, categoría,
# variedad, calidad, unidad, origen, clasificación, min/max/avg price,
# $/Kg, Kg) as the row it displaces, but with its own Fecha (D) and
# Volumen (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31; existing row 31 (and everything below)
# shifts down to row 32 (and below).
$ws.Rows("31:31").Insert()

# Seed the new row 31 with the data that is now in row 32 (the former row
# 31), then overwrite the two cells that actually differ for this new
# observation (Fecha / D and Volumen / J).
$ws.Range("A32:R32").Copy($ws.Range("A31:R31"))

$ws.Range("D31").Value = 44965
$ws.Range("J31").Value = 1200
